$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2792
$ws.Range("I40").Value = 2660
$ws.Range("K40").Value = 2660
$ws.Range("M40").Value = -2485
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H76").Value = 1956192.9
$ws.Range("I76").Value = 2605413.2
$ws.Range("K76").Value = 2605413.2
$ws.Range("M76").Value = -2605098.2
$ws.Range("H79").Value = 1956192.9
$ws.Range("I79").Value = 2605413.2
$ws.Range("K79").Value = 2605413.2
$ws.Range("M79").Value = -2604321.2
$ws.Range("H94").Value = 2627.375
$ws.Range("I94").Value = 2288.4285
$ws.Range("K94").Value = 2288.4285
$ws.Range("M94").Value = -1837.4285
$ws.Range("H106").Value = 3693.077
$ws.Range("I106").Value = 2444.2222
$ws.Range("K106").Value = 2444.2222
$ws.Range("M106").Value = -1813.2222
$ws.Range("H137").Value = 1893.8096
$ws.Range("I137").Value = 1223.5
$ws.Range("J137").Value = 2306.3076
$ws.Range("K137").Value = 3670.5
$ws.Range("L137").Value = 6918.9228
$ws.Range("M137").Value = -1120.5
$ws.Range("N137").Value = -12018.9228

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3767.459
$ws.Range("I32").Value = 2384.3044
$ws.Range("J32").Value = 8009.1333
$ws.Range("K32").Value = 2384.3044
$ws.Range("L32").Value = 8009.1333
$ws.Range("M32").Value = -2097.3044
$ws.Range("N32").Value = -8583.133300000001
$ws.Range("H54").Value = 11980
$ws.Range("J54").Value = 11980
$ws.Range("L54").Value = 11980
$ws.Range("N54").Value = -13518
$ws.Range("H61").Value = 2851.35
$ws.Range("I61").Value = 2246.375
$ws.Range("J61").Value = 5271.25
$ws.Range("K61").Value = 2246.375
$ws.Range("L61").Value = 5271.25
$ws.Range("M61").Value = -2034.375
$ws.Range("N61").Value = -5695.25
$ws.Range("H97").Value = 845.8570999999999
$ws.Range("I97").Value = 486.83334
$ws.Range("K97").Value = 486.83334
$ws.Range("M97").Value = 9.166659999999979
$ws.Range("H102").Value = 1050
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H110").Value = 1462.72
$ws.Range("I110").Value = 1101.4706
$ws.Range("J110").Value = 2230.375
$ws.Range("K110").Value = 1101.4706
$ws.Range("L110").Value = 2230.375
$ws.Range("M110").Value = 943.5293999999999
$ws.Range("N110").Value = -6320.375
$ws.Range("H136").Value = 2851.35
$ws.Range("I136").Value = 2246.375
$ws.Range("J136").Value = 5271.25
$ws.Range("K136").Value = 6739.125
$ws.Range("L136").Value = 15813.75
$ws.Range("M136").Value = -4189.125
$ws.Range("N136").Value = -20913.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 430.66666
$ws.Range("I22").Value = 430.66666
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 430.66666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -257.66666
$ws.Range("N22").ClearContents()
$ws.Range("H37").Value = 28627.5
$ws.Range("I37").Value = 226
$ws.Range("K37").Value = 226
$ws.Range("M37").Value = -89
$ws.Range("H134").Value = 10020.634
$ws.Range("I134").Value = 10795.833
$ws.Range("J134").Value = 6919.8335
$ws.Range("K134").Value = 32387.499
$ws.Range("L134").Value = 20759.5005
$ws.Range("M134").Value = -29852.499
$ws.Range("N134").Value = -25829.5005

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2071545.9
$ws.Range("I58").Value = 2899478
$ws.Range("K58").Value = 2899478
$ws.Range("M58").Value = -2899275
$ws.Range("H94").Value = 854.61536
$ws.Range("I94").Value = 664.8333
$ws.Range("J94").Value = 1017.2857
$ws.Range("K94").Value = 664.8333
$ws.Range("L94").Value = 1017.2857
$ws.Range("M94").Value = -213.8333
$ws.Range("N94").Value = -1919.2857
$ws.Range("H132").Value = 2035.6666
$ws.Range("I132").Value = 1500.3636
$ws.Range("K132").Value = 4501.0908
$ws.Range("M132").Value = -1971.0908
$ws.Range("H134").Value = 964.4
$ws.Range("I134").Value = 964.4
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2893.2
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -358.1999999999998
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 2071545.9
$ws.Range("I136").Value = 2899478
$ws.Range("K136").Value = 8698434
$ws.Range("M136").Value = -8695884

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10522.6875
$ws.Range("I4").Value = 10522.6875
$ws.Range("K4").Value = 31568.0625
$ws.Range("M4").Value = -31456.0625
$ws.Range("H5").Value = 760.9
$ws.Range("J5").Value = 1001.6667
$ws.Range("L5").Value = 3005.0001
$ws.Range("N5").Value = -3229.0001
$ws.Range("H122").Value = 744.25
$ws.Range("I122").Value = 552
$ws.Range("J122").Value = 1064.6666
$ws.Range("K122").Value = 4968
$ws.Range("L122").Value = 9581.999400000001
$ws.Range("M122").Value = -2518
$ws.Range("N122").Value = -14481.9994
$ws.Range("H135").Value = 760.9
$ws.Range("J135").Value = 1001.6667
$ws.Range("L135").Value = 9015.0003
$ws.Range("N135").Value = -14085.0003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3011
$ws.Range("J97").Value = 3011
$ws.Range("L97").Value = 3011
$ws.Range("N97").Value = -4003
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2540.5557
$ws.Range("I122").Value = 2329.5
$ws.Range("K122").Value = 6988.5
$ws.Range("M122").Value = -4538.5
$ws.Range("H132").Value = 9621384
$ws.Range("J132").Value = 7999.3335
$ws.Range("L132").Value = 23998.0005
$ws.Range("N132").Value = -29058.0005

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 529.8570999999999
$ws.Range("I55").Value = 486.33334
$ws.Range("K55").Value = 486.33334
$ws.Range("M55").Value = -313.33334
$ws.Range("H61").Value = 2256.2
$ws.Range("I61").Value = 2188.0667
$ws.Range("K61").Value = 2188.0667
$ws.Range("M61").Value = -1986.0667
$ws.Range("H82").Value = 3664.2222
$ws.Range("I82").Value = 1249.75
$ws.Range("K82").Value = 1249.75
$ws.Range("M82").Value = -888.75
$ws.Range("H85").Value = 3664.2222
$ws.Range("I85").Value = 1249.75
$ws.Range("K85").Value = 1249.75
$ws.Range("M85").Value = -1.75
$ws.Range("H93").Value = 494.9091
$ws.Range("I93").Value = 345
$ws.Range("J93").Value = 1994
$ws.Range("K93").Value = 345
$ws.Range("L93").Value = 1994
$ws.Range("M93").Value = 903
$ws.Range("N93").Value = -4490
$ws.Range("H94").Value = 47526.332
$ws.Range("J94").Value = 47526.332
$ws.Range("L94").Value = 47526.332
$ws.Range("N94").Value = -48878.332
$ws.Range("H95").Value = 94000
$ws.Range("J95").Value = 94000
$ws.Range("L95").Value = 94000
$ws.Range("N95").Value = -99492
$ws.Range("H113").Value = 2256.2
$ws.Range("I113").Value = 2188.0667
$ws.Range("K113").Value = 2188.0667
$ws.Range("M113").Value = -18.06669999999986
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 53531
$ws.Range("J38").Value = 53531
$ws.Range("L38").Value = 53531
$ws.Range("N38").Value = -54477
$ws.Range("H54").Value = 46038.5
$ws.Range("J54").Value = 46038.5
$ws.Range("L54").Value = 46038.5
$ws.Range("N54").Value = -47078.5
$ws.Range("H132").Value = 8235.639999999999
$ws.Range("I132").Value = 2949.5
$ws.Range("J132").Value = 8695.305
$ws.Range("K132").Value = 8848.5
$ws.Range("L132").Value = 26085.915
$ws.Range("M132").Value = -6318.5
$ws.Range("N132").Value = -31145.915
